$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 824.75
$ws.Range("I2").Value = 849.5
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 849.5
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = -736.5
$ws.Range("N2").Value = -1026

$ws.Range("H21").Value = 200017
$ws.Range("I21").Value = 200017
$ws.Range("K21").Value = 200017
$ws.Range("M21").Value = -199549

$ws.Range("H23").Value = 200017
$ws.Range("I23").Value = 200017
$ws.Range("K23").Value = 200017
$ws.Range("M23").Value = -199783

$ws.Range("H62").Value = 5685576
$ws.Range("I62").Value = 11365881
$ws.Range("J62").Value = 5271.273
$ws.Range("K62").Value = 11365881
$ws.Range("L62").Value = 5271.273
$ws.Range("M62").Value = -11365257
$ws.Range("N62").Value = -6519.273

$ws.Range("H65").Value = 5685576
$ws.Range("I65").Value = 11365881
$ws.Range("J65").Value = 5271.273
$ws.Range("K65").Value = 56829405
$ws.Range("L65").Value = 26356.365
$ws.Range("M65").Value = -56826285
$ws.Range("N65").Value = -32596.365

$ws.Range("H137").Value = 5195.148
$ws.Range("I137").Value = 6087.625
$ws.Range("J137").Value = 4819.3687
$ws.Range("K137").Value = 18262.875
$ws.Range("L137").Value = 14458.1061
$ws.Range("M137").Value = -15712.875
$ws.Range("N137").Value = -19558.1061

$ws.Range("H138").Value = 9112.5
$ws.Range("I138").Value = 6551.778
$ws.Range("J138").Value = 10468.177
$ws.Range("K138").Value = 19655.334
$ws.Range("L138").Value = 31404.531
$ws.Range("M138").Value = -14515.334
$ws.Range("N138").Value = -41684.531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3034.524
$ws.Range("I32").Value = 2352.9333
$ws.Range("K32").Value = 2352.9333
$ws.Range("M32").Value = -2065.9333

$ws.Range("H74").Value = 280573.2
$ws.Range("I74").Value = 359746.34
$ws.Range("J74").Value = 3467.125
$ws.Range("K74").Value = 359746.34
$ws.Range("L74").Value = 3467.125
$ws.Range("M74").Value = -358872.34
$ws.Range("N74").Value = -5215.125

$ws.Range("H77").Value = 280573.2
$ws.Range("I77").Value = 359746.34
$ws.Range("J77").Value = 3467.125
$ws.Range("K77").Value = 1798731.7
$ws.Range("L77").Value = 17335.625
$ws.Range("M77").Value = -1794363.7
$ws.Range("N77").Value = -26071.625

$ws.Range("H132").Value = 5248.283
$ws.Range("I132").Value = 1691.0834
$ws.Range("J132").Value = 12781.177
$ws.Range("K132").Value = 5073.2502
$ws.Range("L132").Value = 38343.531
$ws.Range("M132").Value = -2543.2502
$ws.Range("N132").Value = -43403.531

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 58333
$ws.Range("J132").Value = 58333
$ws.Range("L132").Value = 58333
$ws.Range("N132").Value = -68453

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 65625
$ws.Range("J137").Value = 65625
$ws.Range("L137").Value = 65625
$ws.Range("N137").Value = -75825

$ws.Range("H138").Value = 250000
$ws.Range("J138").Value = 250000
$ws.Range("L138").Value = 250000
$ws.Range("N138").Value = -260280

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 678
$ws.Range("I10").Value = 678
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 678
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -539
$ws.Range("N10").ClearContents()

$ws.Range("H19").Value = 5000
$ws.Range("J19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("N19").Value = -5340

$ws.Range("H24").Value = 5000
$ws.Range("J24").Value = 5000
$ws.Range("L24").Value = 5000
$ws.Range("N24").Value = -5340

$ws.Range("H31").Value = 4670.55
$ws.Range("J31").Value = 4951
$ws.Range("L31").Value = 4951
$ws.Range("N31").Value = -5541

$ws.Range("H34").Value = 4670.55
$ws.Range("J34").Value = 4951
$ws.Range("L34").Value = 4951
$ws.Range("N34").Value = -5355

$ws.Range("H100").Value = 50258.332
$ws.Range("J100").Value = 50258.332
$ws.Range("L100").Value = 50258.332
$ws.Range("N100").Value = -52422.332

$ws.Range("H140").Value = 76500
$ws.Range("J140").Value = 76500
$ws.Range("L140").Value = 76500
$ws.Range("N140").Value = -86860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 8491
$ws.Range("I14").Value = 8491
$ws.Range("K14").Value = 25473
$ws.Range("M14").Value = -25300

$ws.Range("H68").Value = 2382.125
$ws.Range("J68").Value = 2600
$ws.Range("L68").Value = 7800
$ws.Range("N68").Value = -9422

$ws.Range("H71").Value = 2382.125
$ws.Range("J71").Value = 2600
$ws.Range("L71").Value = 23400
$ws.Range("N71").Value = -31512

$ws.Range("H112").Value = 5841.8184
$ws.Range("I112").Value = 5851.222
$ws.Range("J112").Value = 5799.5
$ws.Range("K112").Value = 17553.666
$ws.Range("L112").Value = 17398.5
$ws.Range("M112").Value = -16445.666
$ws.Range("N112").Value = -19614.5

$ws.Range("H122").Value = 3900193.5
$ws.Range("I122").Value = 20834046
$ws.Range("J122").Value = 29598.543
$ws.Range("K122").Value = 187506414
$ws.Range("L122").Value = 266386.887
$ws.Range("M122").Value = -187503964
$ws.Range("N122").Value = -271286.887

$ws.Range("H134").Value = 3689.4285
$ws.Range("I134").Value = 3689.4285
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11068.2855
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5998.2855
$ws.Range("N134").ClearContents()

$ws.Range("H141").Value = 2000
$ws.Range("I141").Value = 2000
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6000
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -820
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 70657
$ws.Range("J101").Value = 70657
$ws.Range("L101").Value = 70657
$ws.Range("N101").Value = -77147

$ws.Range("H132").Value = 44907.707
$ws.Range("I132").Value = 3093.65
$ws.Range("J132").Value = 253978
$ws.Range("K132").Value = 9280.95
$ws.Range("L132").Value = 761934
$ws.Range("M132").Value = -6750.950000000001
$ws.Range("N132").Value = -766994

$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8658.143
$ws.Range("I132").Value = 6569.3335
$ws.Range("J132").Value = 10224.75
$ws.Range("K132").Value = 19708.0005
$ws.Range("L132").Value = 30674.25
$ws.Range("M132").Value = -17178.0005
$ws.Range("N132").Value = -35734.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 74398.4
$ws.Range("J125").Value = 74398.4
$ws.Range("L125").Value = 74398.4
$ws.Range("N125").Value = -84238.4

$ws.Range("H126").Value = 4809.7
$ws.Range("I126").Value = 4632.1665
$ws.Range("J126").Value = 5076
$ws.Range("K126").Value = 13896.4995
$ws.Range("L126").Value = 15228
$ws.Range("M126").Value = -11426.4995
$ws.Range("N126").Value = -20168

$ws.Range("H132").Value = 21880.45
$ws.Range("I132").Value = 4177.2954
$ws.Range("J132").Value = 77518.93
$ws.Range("K132").Value = 12531.8862
$ws.Range("L132").Value = 232556.79
$ws.Range("M132").Value = -10001.8862
$ws.Range("N132").Value = -237616.79
